# Update Name of Algo
# Apply updated KNN-imputed values to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = -13.1
$ws.Range("C7").Value = -13.079
$ws.Range("B8").Value = 6.258999999999999
$ws.Range("A12").Value = -21.531
$ws.Range("B12").Value = 6.842000000000001
$ws.Range("B14").Value = 6.622
$ws.Range("C19").Value = -12.955
$ws.Range("D19").Value = -7.722
$ws.Range("C21").Value = -13.079
$ws.Range("B22").Value = 6.593000000000001
$ws.Range("C24").Value = -12.638
